$wb = $excel.ActiveWorkbook

# Rename the sheets (Sheet1 -> CONNECTIVITY, Sheet2 -> COORDINATES, Sheet3 -> FREE NODES)
$wb.Worksheets.Item(1).Name = "CONNECTIVITY"
$wb.Worksheets.Item(2).Name = "COORDINATES"
$wb.Worksheets.Item(3).Name = "FREE NODES"

$wsConn = $wb.Worksheets.Item("CONNECTIVITY")
$wsCoord = $wb.Worksheets.Item("COORDINATES")
$wsFree = $wb.Worksheets.Item("FREE NODES")

# Swap the restraint flags in column E of CONNECTIVITY:
# rows 5-8 become free (1), rows 9-12 become restrained (0)
$wsConn.Range("E5:E8").Value = 1
$wsConn.Range("E9:E12").Value = 0

# Update the remembered selection on FREE NODES sheet
$wsFree.Activate()
$wsFree.Range("H17").Select()

# Update the remembered selection on COORDINATES sheet (no longer the active tab)
$wsCoord.Activate()
$wsCoord.Range("D9").Select()

# CONNECTIVITY becomes the active/selected tab with E13 selected
$wsConn.Activate()
$wsConn.Range("E13").Select()
